$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts every existing column
# (A..W) right by one (becoming B..X) together with styles/merged cells.
$ws.Columns("A").Insert()

# Populate the new "Match ID" column.
$ws.Range("A2").Value = "Match ID"
$ws.Range("A4:A19").Value = 2

# New column gets the bold / borderless header style (matches the new
# cellXfs entry: fontId=1, borderId=0, applyFont=1).
$ws.Range("A2:A18").Font.Bold = $true

# Restore row 19's explicit custom height (writing into A19 above nudges
# the autofit height; put the original 11.4pt back).
$ws.Rows("19").RowHeight = 11.4

# Selection moves to the new Match ID column's data cells.
$ws.Range("A2:A18").Select()
